# Rename the inline logo pictures living in the document's headers/footers.
#
#   footer (PearsonLogo, docPr id=1)      image2.png -> image1.png
#   footer (PearsonLogo, docPr id=2)      image2.png -> image1.png
#   header (BTec_Logo-Orange, docPr id=3) image1.jpg -> image2.jpg
#
# InlineShape.Name is a write-only-ish COM property (it does not reflect
# the picture's existing <wp:docPr name="..."> until this session sets
# it), so shapes are identified by their stable AlternativeText/description
# instead of by their current Name. We walk every section's headers &
# footers so the script does not depend on a particular story-range index.

$d = $word.ActiveDocument

function Rename-LogoShapes {
    param($storyRange, [string]$altText, [string]$newName)

    if ($storyRange -eq $null) { return }
    $shapes = $storyRange.InlineShapes
    if ($shapes -eq $null) { return }

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.AlternativeText -eq $altText) {
            $shp.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            Rename-LogoShapes $hdr.Range "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image1.png"
            Rename-LogoShapes $hdr.Range "BTec_Logo-Orange" "image2.jpg"
        }
    }
    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            Rename-LogoShapes $ftr.Range "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image1.png"
            Rename-LogoShapes $ftr.Range "BTec_Logo-Orange" "image2.jpg"
        }
    }
}
